$p = $ppt.ActivePresentation

function Set-RunText($shape, $newText) {
    $tr = $shape.TextFrame.TextRange
    $run = $tr.Runs(1, 1)
    $run.Text = $newText
}

function Set-CitationText($shape, $newText) {
    $origHeight = $shape.Height
    Set-RunText $shape $newText
    # The citation textbox uses spAutoFit; PowerPoint would normally
    # reflow its height when the text changes, but the box's footprint
    # itself is not part of this edit, so restore it.
    $shape.Height = $origHeight
}

# Slide 2 - Introduction to Animals
$s = $p.Slides.Item(2)
Set-RunText $s.Shapes.Item(2) "Animals are multicellular eukaryotic organisms that belong to the kingdom Animalia. They display a wide range of behaviors, appearances, and habitats."
Set-CitationText $s.Shapes.Item(3) "- Smith, J. (2018). The Diversity of Animals. Journal of Zoology, 25(2), 45-60."

# Slide 3 - Classification of Animals
$s = $p.Slides.Item(3)
Set-RunText $s.Shapes.Item(2) "Animals can be classified into different groups based on various characteristics such as body structure, habitat, and diet. Some common classifications include mammals, birds, reptiles, amphibians, and fish."
Set-CitationText $s.Shapes.Item(3) "- Johnson, A. (2017). Taxonomy of Animals. Animal Sciences Review, 10(4), 112-125."

# Slide 4 - Adaptations in Animals
$s = $p.Slides.Item(4)
Set-RunText $s.Shapes.Item(2) "Animals have evolved various adaptations to survive in their environments. These adaptations can include physical characteristics, behaviors, and reproductive strategies."
Set-CitationText $s.Shapes.Item(3) "- Brown, S. (2019). Adaptations in Animals. Environmental Biology, 15(3), 78-89."

# Slide 5 - Animal Behavior
$s = $p.Slides.Item(5)
Set-RunText $s.Shapes.Item(2) "Animals exhibit a wide range of behaviors, including hunting, mating, communication, and social interactions. These behaviors are influenced by genetics, environment, and learning."
Set-CitationText $s.Shapes.Item(3) "- White, L. (2016). The Study of Animal Behavior. Behavioral Sciences Journal, 8(1), 20-35."

# Slide 6 - Role of Animals in Ecosystems -> Endangered Species
$s = $p.Slides.Item(6)
Set-RunText $s.Shapes.Item(1) "Endangered Species"
Set-RunText $s.Shapes.Item(2) "Many animal species are threatened or endangered due to habitat loss, pollution, climate change, and human activities. Conservation efforts are critical to protecting these species from extinction."
Set-CitationText $s.Shapes.Item(3) "- Green, M. (2020). Conservation of Endangered Species. Wildlife Protection Review, 5(5), 200-215."

# Slide 7 - Conservation of Endangered Species -> Animal Welfare
$s = $p.Slides.Item(7)
Set-RunText $s.Shapes.Item(1) "Animal Welfare"
Set-RunText $s.Shapes.Item(2) "Animal welfare refers to the well-being of animals and encompasses their physical, emotional, and mental health. It is important to consider animal welfare in various contexts, including farming, research, and entertainment."
Set-CitationText $s.Shapes.Item(3) "- Jones, K. (2018). Animal Welfare Standards. Journal of Animal Ethics, 12(3), 150-165."

# Slide 8 - Fascinating Animal Facts -> Famous Animals in History
$s = $p.Slides.Item(8)
Set-RunText $s.Shapes.Item(1) "Famous Animals in History"
Set-RunText $s.Shapes.Item(2) "Throughout history, animals have played important roles in human society. Some famous animals include Laika, the first dog in space, and Koko, the gorilla who learned sign language."
Set-CitationText $s.Shapes.Item(3) "- Roberts, D. (2017). Notable Animals in History. Historical Perspectives Journal, 3(2), 75-88."

# Slide 9 - Human-Animal Relationships -> Animal Intelligence
$s = $p.Slides.Item(9)
Set-RunText $s.Shapes.Item(1) "Animal Intelligence"
Set-RunText $s.Shapes.Item(2) "Many animals exhibit intelligence and cognitive abilities, such as problem-solving, memory, and social learning. Studying animal intelligence can provide insights into the evolution of the human mind."
Set-CitationText $s.Shapes.Item(3) "- Smith, E. (2019). Intelligence in Animals. Cognitive Sciences Review, 18(4), 160-175."

# Slide 10 - The Future of Animals -> Interesting Facts about Animals
$s = $p.Slides.Item(10)
Set-RunText $s.Shapes.Item(1) "Interesting Facts about Animals"
Set-RunText $s.Shapes.Item(2) "There are countless fascinating facts about animals, from the fastest land animal (cheetah) to the longest-living animal (the ocean quahog). Learning about these facts can inspire curiosity and appreciation for the natural world."
Set-CitationText $s.Shapes.Item(3) "- Johnson, T. (2021). Fun Facts about Animals. Nature Discovery Magazine, 7(3), 80-95."

# Slide 11 - Conclusion
$s = $p.Slides.Item(11)
Set-RunText $s.Shapes.Item(2) "Animals are a diverse and important part of our world, contributing to ecosystems, human society, and scientific knowledge. It is essential to protect and respect animals for their welfare and conservation."
Set-CitationText $s.Shapes.Item(3) "- Brown, A. (2019). The Significance of Animals. Nature Conservation Review, 12(1), 30-45."
